# Update the date heading and the 25 two-digit multiplication problems
# in the practice-sheet table. Each Find/Replace targets the full cell
# text (including the trailing "=") so only the intended <w:t> run is
# touched, not a sub-string of another problem.
#
# Note on ordering: "14×16=" is being changed to "96×24=", and the
# document separately already contains an original cell "96×24=" that
# must become "93×83=". Because Find.Execute's ReplaceAll (last arg = 2)
# rewrites every occurrence present in the document at the moment it
# runs, the "96×24=" -> "93×83=" replacement is performed BEFORE the
# "14×16=" -> "96×24=" replacement. That way the newly created "96×24="
# text (from 14×16=) is never caught by the earlier rule, and the
# original "96×24=" cell is safely turned into "93×83=" first.

$d = $word.ActiveDocument

# Header date
$d.Content.Find.Execute("2025-01-17 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-01-18 Saturday", 2)

# Row 1
$d.Content.Find.Execute("51×39=", $true, $false, $false, $false, $false, $true, 1, $false, "53×24=", 2)
$d.Content.Find.Execute("36×47=", $true, $false, $false, $false, $false, $true, 1, $false, "25×73=", 2)
$d.Content.Find.Execute("82×31=", $true, $false, $false, $false, $false, $true, 1, $false, "67×78=", 2)
$d.Content.Find.Execute("98×40=", $true, $false, $false, $false, $false, $true, 1, $false, "81×89=", 2)
$d.Content.Find.Execute("60×86=", $true, $false, $false, $false, $false, $true, 1, $false, "74×92=", 2)

# Row 2
$d.Content.Find.Execute("41×46=", $true, $false, $false, $false, $false, $true, 1, $false, "84×51=", 2)
$d.Content.Find.Execute("26×72=", $true, $false, $false, $false, $false, $true, 1, $false, "41×89=", 2)
$d.Content.Find.Execute("13×65=", $true, $false, $false, $false, $false, $true, 1, $false, "85×51=", 2)
$d.Content.Find.Execute("41×25=", $true, $false, $false, $false, $false, $true, 1, $false, "20×42=", 2)
# "96×24=" (row3 col2's original text) is replaced here, ahead of the
# "14×16=" rule below, to avoid colliding with the text it produces.
$d.Content.Find.Execute("96×24=", $true, $false, $false, $false, $false, $true, 1, $false, "93×83=", 2)
$d.Content.Find.Execute("14×16=", $true, $false, $false, $false, $false, $true, 1, $false, "96×24=", 2)

# Row 3
$d.Content.Find.Execute("27×30=", $true, $false, $false, $false, $false, $true, 1, $false, "63×64=", 2)
$d.Content.Find.Execute("96×91=", $true, $false, $false, $false, $false, $true, 1, $false, "47×88=", 2)
$d.Content.Find.Execute("46×68=", $true, $false, $false, $false, $false, $true, 1, $false, "32×15=", 2)
$d.Content.Find.Execute("98×92=", $true, $false, $false, $false, $false, $true, 1, $false, "25×37=", 2)

# Row 4
$d.Content.Find.Execute("99×26=", $true, $false, $false, $false, $false, $true, 1, $false, "66×55=", 2)
$d.Content.Find.Execute("73×32=", $true, $false, $false, $false, $false, $true, 1, $false, "19×54=", 2)
$d.Content.Find.Execute("40×89=", $true, $false, $false, $false, $false, $true, 1, $false, "99×88=", 2)
$d.Content.Find.Execute("19×67=", $true, $false, $false, $false, $false, $true, 1, $false, "72×90=", 2)
$d.Content.Find.Execute("14×81=", $true, $false, $false, $false, $false, $true, 1, $false, "70×53=", 2)

# Row 5
$d.Content.Find.Execute("47×96=", $true, $false, $false, $false, $false, $true, 1, $false, "36×49=", 2)
$d.Content.Find.Execute("57×63=", $true, $false, $false, $false, $false, $true, 1, $false, "24×84=", 2)
$d.Content.Find.Execute("39×68=", $true, $false, $false, $false, $false, $true, 1, $false, "26×59=", 2)
$d.Content.Find.Execute("41×74=", $true, $false, $false, $false, $false, $true, 1, $false, "85×81=", 2)
$d.Content.Find.Execute("24×55=", $true, $false, $false, $false, $false, $true, 1, $false, "70×61=", 2)
